$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H57").Value = 49999
$ws.Range("J57").Value = 49999
$ws.Range("L57").Value = 149997
$ws.Range("N57").Value = -150995

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4601.3335
$ws.Range("I64").Value = 4601.3335
$ws.Range("K64").Value = 4601.3335
$ws.Range("M64").Value = -4353.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4601.3335
$ws.Range("I67").Value = 4601.3335
$ws.Range("K67").Value = 4601.3335
$ws.Range("M67").Value = -3743.3335

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 10178.375
$ws.Range("I129").Value = 1998.5
$ws.Range("K129").Value = 5995.5
$ws.Range("M129").Value = -995.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 838
$ws.Range("I97").Value = 456.875
$ws.Range("J97").Value = 1752.7
$ws.Range("K97").Value = 456.875
$ws.Range("L97").Value = 1752.7
$ws.Range("M97").Value = 39.125
$ws.Range("N97").Value = -2744.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 166208
$ws.Range("I110").Value = 181940
$ws.Range("K110").Value = 181940
$ws.Range("M110").Value = -179895

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1499.1666
$ws.Range("I5").Value = 1416.3334
$ws.Range("J5").Value = 1582
$ws.Range("K5").Value = 1416.3334
$ws.Range("L5").Value = 1582
$ws.Range("M5").Value = -1303.3334
$ws.Range("N5").Value = -1808

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 203.6
$ws.Range("I22").Value = 243.66667
$ws.Range("J22").Value = 143.5
$ws.Range("K22").Value = 243.66667
$ws.Range("L22").Value = 143.5
$ws.Range("M22").Value = -70.66667000000001
$ws.Range("N22").Value = -489.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 7925.5894
$ws.Range("I107").Value = 9019.152
$ws.Range("J107").Value = 2895.2
$ws.Range("K107").Value = 9019.152
$ws.Range("L107").Value = 2895.2
$ws.Range("M107").Value = -7099.152
$ws.Range("N107").Value = -6735.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3077.9607
$ws.Range("I134").Value = 2763.9268
$ws.Range("K134").Value = 8291.7804
$ws.Range("M134").Value = -5756.7804

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2327.4167
$ws.Range("I31").Value = 1956.0588
$ws.Range("J31").Value = 2579.94
$ws.Range("K31").Value = 1956.0588
$ws.Range("L31").Value = 2579.94
$ws.Range("M31").Value = -1661.0588
$ws.Range("N31").Value = -3169.94

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2327.4167
$ws.Range("I34").Value = 1956.0588
$ws.Range("J34").Value = 2579.94
$ws.Range("K34").Value = 1956.0588
$ws.Range("L34").Value = 2579.94
$ws.Range("M34").Value = -1754.0588
$ws.Range("N34").Value = -2983.94

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 45497880
$ws.Range("I107").Value = 55606830
$ws.Range("J107").Value = 7596.5
$ws.Range("K107").Value = 55606830
$ws.Range("L107").Value = 7596.5
$ws.Range("M107").Value = -55604910
$ws.Range("N107").Value = -11436.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H139").Value = 84299.8
$ws.Range("J139").Value = 90374.75
$ws.Range("L139").Value = 90374.75
$ws.Range("N139").Value = -100654.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1312.375
$ws.Range("I109").Value = 1312.375
$ws.Range("K109").Value = 3937.125
$ws.Range("M109").Value = -2897.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1402.4286
$ws.Range("I113").Value = 608.1667
$ws.Range("J113").Value = 1720.1333
$ws.Range("K113").Value = 1824.5001
$ws.Range("L113").Value = 5160.3999
$ws.Range("M113").Value = 345.4999
$ws.Range("N113").Value = -9500.3999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3056.5217
$ws.Range("I80").Value = 3059.9333
$ws.Range("J80").Value = 3050.125
$ws.Range("K80").Value = 3059.9333
$ws.Range("L80").Value = 3050.125
$ws.Range("M80").Value = -2061.9333
$ws.Range("N80").Value = -5046.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3056.5217
$ws.Range("I83").Value = 3059.9333
$ws.Range("J83").Value = 3050.125
$ws.Range("K83").Value = 15299.6665
$ws.Range("L83").Value = 15250.625
$ws.Range("M83").Value = -10307.6665
$ws.Range("N83").Value = -25234.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 43688.8
$ws.Range("I102").Value = 58780.883
$ws.Range("K102").Value = 58780.883
$ws.Range("M102").Value = -57158.883

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2538.348
$ws.Range("I132").Value = 2740.4106
$ws.Range("K132").Value = 8221.231800000001
$ws.Range("M132").Value = -5691.231800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3279.0476
$ws.Range("I7").Value = 3258
$ws.Range("J7").Value = 3298.182
$ws.Range("K7").Value = 3258
$ws.Range("L7").Value = 3298.182
$ws.Range("M7").Value = -3146
$ws.Range("N7").Value = -3522.182

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2875.25
$ws.Range("J68").Value = 6198.8
$ws.Range("L68").Value = 6198.8
$ws.Range("N68").Value = -7696.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2875.25
$ws.Range("J71").Value = 6198.8
$ws.Range("L71").Value = 30994
$ws.Range("N71").Value = -38482

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H107").Value = 5838.5
$ws.Range("I107").Value = 5838.5
$ws.Range("K107").Value = 5838.5
$ws.Range("M107").Value = -3918.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3279.0476
$ws.Range("I126").Value = 3258
$ws.Range("J126").Value = 3298.182
$ws.Range("K126").Value = 9774
$ws.Range("L126").Value = 9894.545999999998
$ws.Range("M126").Value = -7304
$ws.Range("N126").Value = -14834.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 13339.82
$ws.Range("I136").Value = 1381.3281
$ws.Range("K136").Value = 4143.9843
$ws.Range("M136").Value = -1593.9843

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4910.1113
$ws.Range("I14").Value = 2565.5
$ws.Range("J14").Value = 9599.333000000001
$ws.Range("K14").Value = 2565.5
$ws.Range("L14").Value = 9599.333000000001
$ws.Range("M14").Value = -2397.5
$ws.Range("N14").Value = -9935.333000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 51499.8
$ws.Range("I54").Value = 9999.75
$ws.Range("K54").Value = 9999.75
$ws.Range("M54").Value = -9479.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 756.75
$ws.Range("I107").Value = 582.5714
$ws.Range("J107").Value = 1000.6
$ws.Range("K107").Value = 1747.7142
$ws.Range("L107").Value = 3001.8
$ws.Range("M107").Value = 172.2857999999999
$ws.Range("N107").Value = -6841.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2721.0442
$ws.Range("I122").Value = 2618.8333
$ws.Range("J122").Value = 3115.2856
$ws.Range("K122").Value = 7856.499899999999
$ws.Range("L122").Value = 9345.856800000001
$ws.Range("M122").Value = -5406.499899999999
$ws.Range("N122").Value = -14245.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1101.7819
$ws.Range("I132").Value = 968.0851
$ws.Range("J132").Value = 1887.25
$ws.Range("K132").Value = 2904.2553
$ws.Range("L132").Value = 5661.75
$ws.Range("M132").Value = -374.2552999999998
$ws.Range("N132").Value = -10721.75
